$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the "_old" / "_new" header suffixes to "_FV2210" / "_FV2304".
#    Columns A-J mirror the older "FV2210" release, columns L-U the newer
#    "FV2304" release (K holds the literal "diff" marker and is untouched).
# ---------------------------------------------------------------------------
$fieldNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $fieldNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($fieldNames[$i] + "_FV2210")
}

for ($i = 0; $i -lt $fieldNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($fieldNames[$i] + "_FV2304")
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into a real Excel Table (ListObject) so the sheet
#    gets structured references + an autofilter, matching the regenerated
#    workbook. The header row already carries manual bold/shaded styling;
#    stash it away and restore it afterwards so creating the table doesn't
#    bake that formatting into a new header-row conditional style.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")

$headerRange.Copy()
$scratch.PasteSpecial(-4122)
$headerRange.Style = "Normal"

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U76"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$scratch.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
